$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D and E columns so numeric-looking strings
# (e.g. "138.00", "1.00") keep their exact textual representation
# rather than being converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.844.21'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.773.80'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '426.36'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +8.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.00'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +12.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.737'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +8.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.152'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000310'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.92'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.50'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +15.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.371.24'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.88'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +10.95%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.815.85'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +7.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +7.23%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +11.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.081.69'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '406.19'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.20'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.84%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '36.63'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.83'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +46.99%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.88'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.41'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.93'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +18.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '708.67'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.16%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +16.83%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.95'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +12.47%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +41.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.149'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.54'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0474'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +8.94%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.143'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +9.74%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +50.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.84'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.10%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0671'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.34'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +9.23%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.317'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +15.49%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.67%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.81'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.80'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.77%  '
